$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1315.2941
$ws.Range("I2").Value = 410.57144
$ws.Range("J2").Value = 1948.6
$ws.Range("K2").Value = 410.57144
$ws.Range("L2").Value = 1948.6
$ws.Range("M2").Value = -297.57144
$ws.Range("N2").Value = -2174.6
$ws.Range("H17").Value = 2050
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H29").Value = 706.8333
$ws.Range("I29").Value = 185.25
$ws.Range("K29").Value = 555.75
$ws.Range("M29").Value = -274.75
$ws.Range("H55").Value = 647.1429000000001
$ws.Range("I55").Value = 407.5
$ws.Range("K55").Value = 407.5
$ws.Range("M55").Value = -193.5
$ws.Range("H70").Value = 3114.818
$ws.Range("I70").Value = 2186.625
$ws.Range("J70").Value = 3645.2144
$ws.Range("K70").Value = 6559.875
$ws.Range("L70").Value = 10935.6432
$ws.Range("M70").Value = -6289.875
$ws.Range("N70").Value = -11475.6432
$ws.Range("H73").Value = 3114.818
$ws.Range("I73").Value = 2186.625
$ws.Range("J73").Value = 3645.2144
$ws.Range("K73").Value = 6559.875
$ws.Range("L73").Value = 10935.6432
$ws.Range("M73").Value = -5623.875
$ws.Range("N73").Value = -12807.6432
$ws.Range("H80").Value = 554.0909
$ws.Range("I80").Value = 191.85715
$ws.Range("J80").Value = 1188
$ws.Range("K80").Value = 575.5714499999999
$ws.Range("L80").Value = 3564
$ws.Range("M80").Value = 422.4285500000001
$ws.Range("N80").Value = -5560
$ws.Range("H83").Value = 554.0909
$ws.Range("I83").Value = 191.85715
$ws.Range("J83").Value = 1188
$ws.Range("K83").Value = 1726.71435
$ws.Range("L83").Value = 10692
$ws.Range("M83").Value = 3265.28565
$ws.Range("N83").Value = -20676
$ws.Range("H86").Value = 25251
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 25251
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 25251
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -27497
$ws.Range("H89").Value = 25251
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 25251
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 126255
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -137487
$ws.Range("H100").Value = 9379.6
$ws.Range("I100").Value = 3966
$ws.Range("J100").Value = 17500
$ws.Range("K100").Value = 3966
$ws.Range("L100").Value = 17500
$ws.Range("M100").Value = -3425
$ws.Range("N100").Value = -18582
$ws.Range("H135").Value = 1198.2222
$ws.Range("I135").Value = 297.83334
$ws.Range("K135").Value = 2680.50006
$ws.Range("M135").Value = -145.5000600000003
$ws.Range("H137").Value = 2242.25
$ws.Range("I137").Value = 1682.0834
$ws.Range("J137").Value = 3922.75
$ws.Range("K137").Value = 5046.2502
$ws.Range("L137").Value = 11768.25
$ws.Range("M137").Value = -2496.2502
$ws.Range("N137").Value = -16868.25
$ws.Range("H138").Value = 2810
$ws.Range("I138").Value = 1372.2
$ws.Range("J138").Value = 3289.2666
$ws.Range("K138").Value = 4116.6
$ws.Range("L138").Value = 9867.799800000001
$ws.Range("M138").Value = 1023.4
$ws.Range("N138").Value = -20147.7998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7782.2
$ws.Range("I32").Value = 7782.2
$ws.Range("K32").Value = 7782.2
$ws.Range("M32").Value = -7495.2
$ws.Range("H88").Value = 7998.5
$ws.Range("J88").Value = 7998.5
$ws.Range("L88").Value = 7998.5
$ws.Range("N88").Value = -8810.5
$ws.Range("H91").Value = 7998.5
$ws.Range("J91").Value = 7998.5
$ws.Range("L91").Value = 7998.5
$ws.Range("N91").Value = -10806.5
$ws.Range("H132").Value = 3489.0476
$ws.Range("I132").Value = 3113.5
$ws.Range("K132").Value = 9340.5
$ws.Range("M132").Value = -6810.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4250
$ws.Range("I86").Value = 4000
$ws.Range("J86").Value = 4500
$ws.Range("K86").Value = 4000
$ws.Range("L86").Value = 4500
$ws.Range("M86").Value = -2877
$ws.Range("N86").Value = -6746
$ws.Range("H89").Value = 4250
$ws.Range("I89").Value = 4000
$ws.Range("J89").Value = 4500
$ws.Range("K89").Value = 20000
$ws.Range("L89").Value = 22500
$ws.Range("M89").Value = -14384
$ws.Range("N89").Value = -33732
$ws.Range("H94").Value = 2195.1428
$ws.Range("I94").Value = 2215.9473
$ws.Range("J94").Value = 1997.5
$ws.Range("K94").Value = 2215.9473
$ws.Range("L94").Value = 1997.5
$ws.Range("M94").Value = -1764.9473
$ws.Range("N94").Value = -2899.5
$ws.Range("H99").Value = 855
$ws.Range("I99").Value = 855
$ws.Range("K99").Value = 855
$ws.Range("M99").Value = 643
$ws.Range("H134").Value = 10220.333
$ws.Range("I134").Value = 12227
$ws.Range("J134").Value = 6207
$ws.Range("K134").Value = 36681
$ws.Range("L134").Value = 18621
$ws.Range("M134").Value = -34146
$ws.Range("N134").Value = -23691

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 25408.2
$ws.Range("J28").Value = 25408.2
$ws.Range("L28").Value = 25408.2
$ws.Range("N28").Value = -25898.2
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H96").Value = 20624
$ws.Range("J96").Value = 20624
$ws.Range("L96").Value = 20624
$ws.Range("N96").Value = -26116
$ws.Range("H132").Value = 1920.6
$ws.Range("I132").Value = 1920.6
$ws.Range("K132").Value = 5761.799999999999
$ws.Range("M132").Value = -3231.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1361.9688
$ws.Range("J4").Value = 2218.182
$ws.Range("L4").Value = 6654.545999999999
$ws.Range("N4").Value = -6878.545999999999
$ws.Range("H75").Value = 500
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 500
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 1500
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -3496
$ws.Range("H78").Value = 500
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 500
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 4500
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -14484
$ws.Range("H92").Value = 499
$ws.Range("J92").Value = 499
$ws.Range("L92").Value = 1497
$ws.Range("N92").Value = -3993
$ws.Range("H110").Value = 2800
$ws.Range("I110").Value = 2800
$ws.Range("K110").Value = 8400
$ws.Range("M110").Value = -4310
$ws.Range("H138").Value = 1583
$ws.Range("I138").Value = 612.3333
$ws.Range("J138").Value = 1999
$ws.Range("K138").Value = 1836.9999
$ws.Range("L138").Value = 5997
$ws.Range("M138").Value = 3303.0001
$ws.Range("N138").Value = -16277

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4817.909
$ws.Range("I80").Value = 2998.5
$ws.Range("J80").Value = 5222.222
$ws.Range("K80").Value = 2998.5
$ws.Range("L80").Value = 5222.222
$ws.Range("M80").Value = -2000.5
$ws.Range("N80").Value = -7218.222
$ws.Range("H83").Value = 4817.909
$ws.Range("I83").Value = 2998.5
$ws.Range("J83").Value = 5222.222
$ws.Range("K83").Value = 14992.5
$ws.Range("L83").Value = 26111.11
$ws.Range("M83").Value = -10000.5
$ws.Range("N83").Value = -36095.11
$ws.Range("H97").Value = 808.86957
$ws.Range("I97").Value = 830.25
$ws.Range("J97").Value = 666.3333
$ws.Range("K97").Value = 830.25
$ws.Range("L97").Value = 666.3333
$ws.Range("M97").Value = -334.25
$ws.Range("N97").Value = -1658.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 845.5
$ws.Range("I22").Value = 600.64703
$ws.Range("J22").Value = 1308
$ws.Range("K22").Value = 600.64703
$ws.Range("L22").Value = 1308
$ws.Range("M22").Value = -305.64703
$ws.Range("N22").Value = -1898
$ws.Range("H27").Value = 845.5
$ws.Range("I27").Value = 600.64703
$ws.Range("J27").Value = 1308
$ws.Range("K27").Value = 600.64703
$ws.Range("L27").Value = 1308
$ws.Range("M27").Value = -493.64703
$ws.Range("N27").Value = -1522
$ws.Range("H46").Value = 1149
$ws.Range("J46").Value = 3998
$ws.Range("L46").Value = 3998
$ws.Range("N46").Value = -4374
$ws.Range("H82").Value = 1507.65
$ws.Range("I82").Value = 1331
$ws.Range("J82").Value = 1835.7142
$ws.Range("K82").Value = 1331
$ws.Range("L82").Value = 1835.7142
$ws.Range("M82").Value = -970
$ws.Range("N82").Value = -2557.7142
$ws.Range("H85").Value = 1507.65
$ws.Range("I85").Value = 1331
$ws.Range("J85").Value = 1835.7142
$ws.Range("K85").Value = 1331
$ws.Range("L85").Value = 1835.7142
$ws.Range("M85").Value = -83
$ws.Range("N85").Value = -4331.7142
$ws.Range("H93").Value = 1000
$ws.Range("I93").Value = 1000
$ws.Range("K93").Value = 1000
$ws.Range("M93").Value = 248
$ws.Range("H132").Value = 12724.735
$ws.Range("I132").Value = 14223.409
$ws.Range("J132").Value = 9977.166999999999
$ws.Range("K132").Value = 42670.227
$ws.Range("L132").Value = 29931.501
$ws.Range("M132").Value = -40140.227
$ws.Range("N132").Value = -34991.501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 650
$ws.Range("I81").Value = 800
$ws.Range("J81").Value = 500
$ws.Range("K81").Value = 1600
$ws.Range("L81").Value = 1000
$ws.Range("M81").Value = -539
$ws.Range("N81").Value = -3122
$ws.Range("H84").Value = 650
$ws.Range("I84").Value = 800
$ws.Range("J84").Value = 500
$ws.Range("K84").Value = 8000
$ws.Range("L84").Value = 5000
$ws.Range("M84").Value = -2696
$ws.Range("N84").Value = -15608
$ws.Range("H104").Value = 13332
$ws.Range("J104").Value = 13332
$ws.Range("L104").Value = 13332
$ws.Range("N104").Value = -20320
